$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift data one column to the left (old B:F -> new A:E) for rows 1-3.
# This drops the old column A (the per-row "GENE id" numbers, which used to
# carry the header-style border/bold/center/top formatting) and slides the
# old column F data into the new column E.
for ($r = 1; $r -le 3; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $f = $ws.Cells.Item($r, 6).Value2

    $ws.Cells.Item($r, 1).Value2 = $b
    $ws.Cells.Item($r, 2).Value2 = $c
    $ws.Cells.Item($r, 3).Value2 = $d
    $ws.Cells.Item($r, 4).Value2 = $e
    $ws.Cells.Item($r, 5).Value2 = $f
}

# The header row (row 1) keeps the bold/border/center/top style across the
# whole A1:E1 span now, so copy that formatting onto the new A1 cell too.
$ws.Range("B1:E1").Copy()
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

# The old column A (row2/row3) carried the header-ish style for the GENE id
# - that column is gone now, so clear the leftover formatting on the new A2/A3.
$ws.Range("A2:A3").ClearFormats()

# Column F no longer holds any data - clear its contents/formatting entirely.
$ws.Range("F1:F3").Clear()
